$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 65 (shifts old rows 65:97 down to 66:98)
$ws.Rows.Item(65).Insert()

# Populate the newly inserted row 65 with the new record's data
$ws.Range("A65").Value = 5
$ws.Range("B65").Value = "Macroferia Regional de Talca"
$ws.Range("C65").Value = "Maule"
$ws.Range("D65").Value = 44917
$ws.Range("D65").NumberFormat = $ws.Range("D66").NumberFormat
$ws.Range("E65").Value = 7
$ws.Range("F65").Value = "Fruta"
$ws.Range("G65").Value = 100101
$ws.Range("H65").Value = "Berries"
$ws.Range("I65").Value = 100101001
$ws.Range("J65").Value = "Arándano (blue)"
$ws.Range("K65").Value = "Sin especificar"
$ws.Range("L65").Value = "Primera"
$ws.Range("M65").Value = 60
$ws.Range("N65").Value = 3000
$ws.Range("O65").Value = 3000
$ws.Range("P65").Value = 3000
$ws.Range("Q65").Value = "$/bandeja 2 kilos"
$ws.Range("R65").Value = "Provincia de Curicó"
$ws.Range("S65").Value = 1500
$ws.Range("T65").Value = 2
